$d = $word.ActiveDocument

# --- Update the date heading paragraph ---
$d.Paragraphs.Item(1).Range.Text = "2023-11-19 Sunday"

# --- Update the practice-problem table ---
# The table has 20 rows; data lives in rows 1, 5, 9, 13, 17 (5 cells each).
# All 25 data cells get new text (row 13's values also shift by one
# position because the original first cell of that row was removed and a
# new cell was appended at the end), so every cell is addressed directly
# by (row, column) position rather than via text search/replace to avoid
# any ambiguity from repeated/overlapping values.
$tbl = $d.Tables.Item(1)

$rowData = @{
    1  = @("69÷8=8, 5", "65÷2=32, 1", "43÷6=7, 1", "29÷6=4, 5", "98÷6=16, 2")
    5  = @("37÷6=6, 1", "22÷6=3, 4", "13÷6=2, 1", "18÷8=2, 2", "52÷4=13, 0")
    9  = @("42÷9=4, 6", "60÷5=12, 0", "20÷3=6, 2", "47÷9=5, 2", "47÷6=7, 5")
    13 = @("58÷6=9, 4", "24÷9=2, 6", "69÷3=23, 0", "34÷2=17, 0", "62÷4=15, 2")
    17 = @("25÷9=2, 7", "74÷7=10, 4", "44÷6=7, 2", "40÷6=6, 4", "21÷7=3, 0")
}

foreach ($rowIndex in $rowData.Keys) {
    $row = $tbl.Rows.Item($rowIndex)
    $values = $rowData[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $row.Cells.Item($col).Range.Text = $values[$col - 1]
    }
}
